$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (price) as text so number-like strings (e.g. "1.930", "7.470")
# are not silently coerced into floating point numbers, which would lose
# trailing zeros / thousands-style dot grouping present in the source data.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "25.834.84"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "1.633.81"
$ws.Range("D5").Value = "214.62"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "0.5016"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "0.06387"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Value = "19.66"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "0.07677"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.244"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.631.69"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "1.858.12"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "0.5424"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "0.0₅7920"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "63.54"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "25.847.73"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "201.74"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("D21").Value = "4.333"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "9.918"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "5.967"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "1.930"
$ws.Range("E25").Value = "  +10.99%  "
$ws.Range("D26").Value = "141.78"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "6.708"
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "0.04997"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").Value = "3.260"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("D33").Value = "3.182"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "1.539"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "1.175.29"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").Value = "0.8915"
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("E38").Value = "  -4.88%  "
$ws.Range("D39").Value = "0.5587"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("D40").Value = "0.01561"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "5.692"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "0.8063"
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("D44").Value = "99.52"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "1.769.95"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").Value = "0.4513"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "54.78"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("D50").Value = "0.05076"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "7.470"
$ws.Range("E51").Value = "  -4.97%  "

# Restore default (unstyled) appearance for column D now that values are set.
$priceRange.Style = "Normal"

